# Team19 poster - slide 1: reference Fig.3/Fig.2 explicitly instead of
# "the pie chart above" / "the accompanying picture".
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Locate the "Text Placeholder 18" shape that holds the results paragraph
# mentioning "the pie chart above" (avoid relying on a hard-coded index).
$target = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $sh = $s.Shapes.Item($i)
    if ($sh.HasTextFrame) {
        if ($sh.TextFrame.TextRange.Text -like "*the pie chart above*") {
            $target = $sh
            break
        }
    }
}

$tr = $target.TextFrame.TextRange

# "...as you can see from the pie chart above, these are..."
#   -> "...as you can see from Fig.3 above, these are..."
$tr.Replace("the pie chart above", "Fig.3 above", 0, $false, $false) | Out-Null

# "...in the accompanying picture, measured..."
#   -> "...in the accompanying picture in Fig.2, measured..."
$tr2 = $target.TextFrame.TextRange
$tr2.Replace("accompanying picture, ", "accompanying picture in Fig.2, ", 0, $false, $false) | Out-Null
